$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# Enter the reported totals for Day 1 (rows 3-5)
$ws.Range("C3").Value = 7075
$ws.Range("C4").Value = 2610
$ws.Range("C5").Value = 2610

# Row 3 height adjusts to fit the new content (37.5 -> 30)
$ws.Rows.Item(3).RowHeight = 30

# Re-establish the merged header cells so they are stored in row order
$mergedRanges = @("B2:C2","B8:C8","B14:C14","B20:C20","B26:C26","B32:C32","B38:C38","B44:C44","B50:C50","B57:C57")
foreach ($r in $mergedRanges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $mergedRanges) {
    $ws.Range($r).Merge()
}

# Scroll the view back to the top and select C6
$ws.Range("C6").Select()
